$d = $word.ActiveDocument

function Replace-InParagraph($paraIndex, $oldText, $newText) {
    $p = $d.Paragraphs.Item($paraIndex)
    $rng = $d.Range($p.Range.Start, $p.Range.End)
    $ok = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 1)
    if (-not $ok) {
        Write-Output ("REPLACE FAILED in paragraph " + $paraIndex)
    }
}

# Paragraph 6, run 1
$oldVal = "Aplicar os conhecimentos adquiridos nas disciplinas obrigatórias e demonstrar a capacidade de articulação das competências inerentes ao profissional de Engenharia Bioquímica através da apresentação de um trabalho acadêmico."
$newVal = "Elaborar e desenvolver, individualmente ou em equipe, um projeto de engenharia, de pesquisa científica ou tecnológica que aplique de forma eficiente o conhecimento técnico e o pensamento crítico, reflexivo, e ético adquiridos durante o processo formativo."
Replace-InParagraph 6 $oldVal $newVal

# Paragraph 7, run 1
$oldVal = "To apply the knowledge acquired in mandatory disciplines and to demonstrate the ability to articulate the skills inherent to the Biochemical Engineering professional through the presentation of an academic paper."
$newVal = "To elaborate and develop, individually or as a team, an engineering project, scientific research, or a business model, and efficiently apply the technical expertise and critical, reflective, and ethical thinking acquired during the educational process."
Replace-InParagraph 7 $oldVal $newVal

# Paragraph 9, run 1
$oldVal = "3380737 - Flávio Teixeira da Silva" + [char]11
$newVal = "Aplicar os conhecimentos adquiridos nas disciplinas obrigatórias e demonstrar a capacidade de articulação das competências inerentes ao profissional de Engenharia Bioquímica através da apresentação de um trabalho acadêmico." + [char]11
Replace-InParagraph 9 $oldVal $newVal

# Paragraph 9, run 2
$oldVal = "8853480 - Tatiane da Franca Silva"
$newVal = "A disciplina consiste no desenvolvimento de um trabalho acadêmico supervisionado por um docente e/ou profissional com título de doutor (orientador), nos seguintes contextos: " + [char]11 + "1. Desenvolvimento de um projeto de Engenharia, podendo abordar a problemática trabalhada na disciplina LOT2062 Solução de Problemas de Engenharia; " + [char]11 + "2. Realização de pesquisa científica ou tecnológica inédita, de caráter teórico ou experimental, abordando temas relevantes na área de Engenharia que demandem atualização e síntese de informações."
Replace-InParagraph 9 $oldVal $newVal

# Paragraph 11, run 1
$oldVal = "Elaborar e desenvolver, individualmente ou em equipe, um projeto de engenharia, de pesquisa científica ou tecnológica que aplique de forma eficiente o conhecimento técnico e o pensamento crítico, reflexivo, e ético adquiridos durante o processo formativo."
$newVal = "Análise da monografia e defesa do trabalho perante uma banca de 3 examinadores, conforme normas específicas."
Replace-InParagraph 11 $oldVal $newVal

# Paragraph 12, run 1
$oldVal = "To elaborate and develop, individually or as a team, an engineering project, scientific research, or a business model, and efficiently apply the technical expertise and critical, reflective, and ethical thinking acquired during the educational process."
$newVal = "To apply the knowledge acquired in mandatory disciplines and to demonstrate the ability to articulate the skills inherent to the Biochemical Engineering professional through the presentation of an academic paper."
Replace-InParagraph 12 $oldVal $newVal

# Paragraph 14, run 1
$oldVal = "A disciplina consiste no desenvolvimento de um trabalho acadêmico supervisionado por um docente e/ou profissional com título de doutor (orientador), nos seguintes contextos: " + [char]11 + "1. Desenvolvimento de um projeto de Engenharia, podendo abordar a problemática trabalhada na disciplina LOT2062 Solução de Problemas de Engenharia; " + [char]11 + "2. Realização de pesquisa científica ou tecnológica inédita, de caráter teórico ou experimental, abordando temas relevantes na área de Engenharia que demandem atualização e síntese de informações."
$newVal = "A nota será individual e atribuída pelos docentes da banca examinadora."
Replace-InParagraph 14 $oldVal $newVal

# Paragraph 17, run 2
$oldVal = "Análise da monografia e defesa do trabalho perante uma banca de 3 examinadores, conforme normas específicas." + [char]11
$newVal = "Reapresentação do trabalho e/ou da monografia para nova avaliação." + [char]11
Replace-InParagraph 17 $oldVal $newVal

# Paragraph 17, run 4
$oldVal = "A nota será individual e atribuída pelos docentes da banca examinadora." + [char]11
$newVal = "A ser definida em função do projeto" + [char]11
Replace-InParagraph 17 $oldVal $newVal

# Paragraph 17, run 6
$oldVal = "Reapresentação do trabalho e/ou da monografia para nova avaliação."
$newVal = "3380737 - Flávio Teixeira da Silva"
Replace-InParagraph 17 $oldVal $newVal

# Paragraph 19, run 1
$oldVal = "A ser definida em função do projeto"
$newVal = "8853480 - Tatiane da Franca Silva"
Replace-InParagraph 19 $oldVal $newVal

Write-Output "DONE"